# Insert a new row above row 43 for the "plan specification" (IAO:0000104)
# entry, pushing the existing rows 43-55 down to 44-56.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 43; everything from the old row 43
# downward shifts down by one.
$ws.Rows(43).Insert()

# The inserted row copies the formatting of the row above it (row 42);
# the new row in the target workbook carries no explicit style, so clear
# any inherited formatting.
$ws.Rows(43).ClearFormats()

# Populate the new row 43 with the "plan specification" entry.
$ws.Cells.Item(43, 1).Value = "IAO:0000104"
$ws.Cells.Item(43, 2).Value = "plan specification"
$ws.Cells.Item(43, 3).Value = "A directive information entity with action specifications and objective specifications as parts, and that may be concretized as a realizable entity that, if realized, is realized in a process in which the bearer tries to achieve the objectives by taking the actions specified."
$ws.Cells.Item(43, 4).Value = "directive information entity"
$ws.Cells.Item(43, 16).Value = "LSR 1"
$ws.Cells.Item(43, 17).Value = "Intervention content and delivery"
$ws.Cells.Item(43, 19).Value = "External"
$ws.Cells.Item(43, 22).Value = "PS"
